$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 135, shifting existing rows 135-144 down to 136-145.
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new weekly data point
# (same template columns as the surrounding rows: Mercado ID/Mercado/Region/
# Codreg/Categoria ID/Categoria/Variedad/Calidad/Origen/Clasificacion).
$ws.Cells.Item(135, 1).Value = 8
$ws.Cells.Item(135, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(135, 3).Value = "Coquimbo"
$ws.Cells.Item(135, 4).Value = 44746
$ws.Cells.Item(135, 5).Value = 4
$ws.Cells.Item(135, 6).Value = 100112001
$ws.Cells.Item(135, 7).Value = "Berenjena"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 700
$ws.Cells.Item(135, 11).Value = 9000
$ws.Cells.Item(135, 12).Value = 10000
$ws.Cells.Item(135, 13).Value = 9500
$ws.Cells.Item(135, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(135, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(135, 16).Value = 190
$ws.Cells.Item(135, 17).Value = 50
$ws.Cells.Item(135, 18).Value = "Hortaliza"
